# Update Tarefa2 Criação 4
# Adds the "Método Recursivo" (recursive GCD) worked example under the
# existing "Teste de Mesa - Fração" table, and tidies up the G5 cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Método Recursivo" section -------------------------------------
# Shared-string table order matters: write "a" (and friends) before
# "Método Recursivo" so new <si> entries land in the same order the
# original commit produced them. Only cells that actually hold data are
# touched, so no incidental blank/styled cells get written.

$ws.Range("B13").Value = "a"
$ws.Range("B12").Value = "Método Recursivo"
$ws.Range("C13").Value = "b"
$ws.Range("D13").Value = "mdc"
$ws.Range("E13").Value = "resto"

$ws.Range("B14").Value = 4
$ws.Range("C14").Value = 8
$ws.Range("E14").Value = 4

$ws.Range("B15").Value = 8
$ws.Range("C15").Value = 4
$ws.Range("E15").Value = 0

$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 4

# Match the existing table formatting (centered horizontally & vertically)
# — applied cell-by-cell so we don't touch neighbouring blank cells.
$newCells = "B12","B13","C13","D13","E13","B14","C14","E14","B15","C15","E15","B16","C16"
foreach ($addr in $newCells) {
  $ws.Range($addr).HorizontalAlignment = -4108
  $ws.Range($addr).VerticalAlignment = -4108
}

# --- Tidy up G5 (drop the stray number-format flag, keep centered align) -

$ws.Range("G5").HorizontalAlignment = -4108
$ws.Range("G5").VerticalAlignment = -4108

# --- View state: zoom out a bit and move the selection -------------------

$excel.ActiveWindow.Zoom = 91
[void]$ws.Range("C16").Select()

Write-Output "done"
